$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startYear = 1987
$endYear = 2024

$row = 2
for ($year = $startYear; $year -le $endYear; $year++) {
    $label = "$($year)Q4"
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $label

    # Match the header cell's formatting (plain text style, no custom
    # date number format) instead of the inherited date-serial style.
    $ws.Range("A1").Copy()
    $cell.PasteSpecial(-4122)

    $row = $row + 1
}

$excel.CutCopyMode = $false

# The date-serial custom number format is no longer referenced by any
# cell now that column A holds text labels; drop it from the style
# table (mirrors what Excel does when the last cell using a custom
# format is reformatted).
$wb.DeleteNumberFormat("YYYY-MM-DD HH:MM:SS")
